# This workbook's data rows (2-44) were re-shuffled: every row keeps its
# "Mercado/Producto/Categoría/Variedad" columns (A,B,C,E,F,G,H,I,J,K — they
# are identical for every row anyway) but the per-record columns
# (D=Fecha, L=Calidad, M=Volumen, N=Precio mínimo, O=Precio máximo,
#  P=Precio promedio ponderado, Q=Unidad de comercialización, R=Origen,
#  S=Precio $/Kg, T=Kg/unidad) move to a different row according to a
# fixed permutation. Capture a full snapshot first (so source values
# aren't clobbered while writing), then write each destination row from
# its mapped source row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns involved in the per-row payload that moves: D, L, M, N, O, P, Q, R, S, T
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

$firstRow = 2
$lastRow = 44

# 1) Snapshot every source row's payload values before writing anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowvals = @{}
    foreach ($c in $cols) {
        $rowvals[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowvals
}

# 2) new row -> old row it inherits its payload from.
$mapping = @{
    2 = 10;  3 = 17;  4 = 18;  5 = 31;  6 = 24;  7 = 33;  8 = 8;   9 = 4;   10 = 23;
    11 = 34; 12 = 35; 13 = 43; 14 = 44; 15 = 3;  16 = 25; 17 = 26; 18 = 20; 19 = 15;
    20 = 16; 21 = 42; 22 = 21; 23 = 36; 24 = 2;  25 = 39; 26 = 40; 27 = 13; 28 = 14;
    29 = 27; 30 = 38; 31 = 32; 32 = 7;  33 = 9;  34 = 5;  35 = 41; 36 = 6;  37 = 19;
    38 = 28; 39 = 29; 40 = 30; 41 = 37; 42 = 22; 43 = 11; 44 = 12
}

# 3) Write each destination row's payload from the snapshot of its source row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $mapping[$r]
    $rowvals = $snapshot[$src]
    foreach ($c in $cols) {
        $ws.Cells.Item($r, $c).Value2 = $rowvals[$c]
    }
}
